$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new schedule row (第82期) below the existing data.
$ws.Range("A32").Value = "10/31"
$ws.Range("C32").Value = "第82期 第二代星途"
$ws.Range("B32").Value = "12/26"

# Match the saved selection state recorded in the workbook.
[void]$ws.Range("C27").Select()
